$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
$newDate = Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0
$newDate = $newDate.Date
$ws.Range("C2:C5").Value = $newDate
